# Auto-generated script applying the Brynhildr_Profits cell-value updates
# scraped from the scheduled runner diff (currentAveragePrice / Leve profit columns).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 131.125
$ws.Range("I33").Value = 144
$ws.Range("J33").Value = 92.5
$ws.Range("K33").Value = 144
$ws.Range("L33").Value = 92.5
$ws.Range("M33").Value = 85
$ws.Range("N33").Value = -550.5
$ws.Range("H106").Value = 12374.75
$ws.Range("I106").Value = 7799.6
$ws.Range("K106").Value = 7799.6
$ws.Range("M106").Value = -7168.6
$ws.Range("H113").Value = 5109.375
$ws.Range("I113").Value = 5235
$ws.Range("J113").Value = 4230
$ws.Range("K113").Value = 5235
$ws.Range("L113").Value = 4230
$ws.Range("M113").Value = -1981
$ws.Range("N113").Value = -10738
$ws.Range("H116").Value = 24848.46
$ws.Range("I116").Value = 12379
$ws.Range("J116").Value = 44799.6
$ws.Range("K116").Value = 12379
$ws.Range("L116").Value = 44799.6
$ws.Range("M116").Value = -8937
$ws.Range("N116").Value = -51683.6
$ws.Range("H118").Value = 1250841.8
$ws.Range("I118").Value = 1666789
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 5000367
$ws.Range("L118").Value = 9000
$ws.Range("M118").Value = -4998710
$ws.Range("N118").Value = -12314
$ws.Range("H132").Value = 1322.3478
$ws.Range("I132").Value = 1306.1904
$ws.Range("J132").Value = 1492
$ws.Range("K132").Value = 3918.5712
$ws.Range("L132").Value = 4476
$ws.Range("M132").Value = -1388.5712
$ws.Range("N132").Value = -9536
$ws.Range("H137").Value = 1707.8788
$ws.Range("I137").Value = 1098.6072
$ws.Range("K137").Value = 3295.8216
$ws.Range("M137").Value = -745.8215999999998
$ws.Range("H138").Value = 5095.932
$ws.Range("J138").Value = 6028
$ws.Range("L138").Value = 18084
$ws.Range("N138").Value = -28364

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 7226.3335
$ws.Range("I63").Value = 4005
$ws.Range("K63").Value = 4005
$ws.Range("M63").Value = -3319
$ws.Range("H66").Value = 7226.3335
$ws.Range("I66").Value = 4005
$ws.Range("K66").Value = 20025
$ws.Range("M66").Value = -16593
$ws.Range("H74").Value = 3367.3513
$ws.Range("I74").Value = 2670.2
$ws.Range("K74").Value = 2670.2
$ws.Range("M74").Value = -1796.2
$ws.Range("H77").Value = 3367.3513
$ws.Range("I77").Value = 2670.2
$ws.Range("K77").Value = 13351
$ws.Range("M77").Value = -8983
$ws.Range("H97").Value = 740.9583
$ws.Range("I97").Value = 762.6087
$ws.Range("K97").Value = 762.6087
$ws.Range("M97").Value = -266.6087
$ws.Range("H102").Value = 913.8261
$ws.Range("I102").Value = 923.5454999999999
$ws.Range("K102").Value = 923.5454999999999
$ws.Range("M102").Value = 698.4545000000001
$ws.Range("H110").Value = 1649.44
$ws.Range("I110").Value = 1539
$ws.Range("K110").Value = 1539
$ws.Range("M110").Value = 506
$ws.Range("H122").Value = 2649.4546
$ws.Range("I122").Value = 2530.4
$ws.Range("J122").Value = 2904.5715
$ws.Range("K122").Value = 7591.200000000001
$ws.Range("L122").Value = 8713.7145
$ws.Range("M122").Value = -5141.200000000001
$ws.Range("N122").Value = -13613.7145

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1370.3334
$ws.Range("I37").Value = 1180.75
$ws.Range("J37").Value = 1749.5
$ws.Range("K37").Value = 1180.75
$ws.Range("L37").Value = 1749.5
$ws.Range("M37").Value = -1043.75
$ws.Range("N37").Value = -2023.5
$ws.Range("H105").Value = 3464.8823
$ws.Range("I105").Value = 3407.111
$ws.Range("J105").Value = 3529.875
$ws.Range("K105").Value = 3407.111
$ws.Range("L105").Value = 3529.875
$ws.Range("M105").Value = -1660.111
$ws.Range("N105").Value = -7023.875
$ws.Range("H107").Value = 1339.4286
$ws.Range("I107").Value = 1307.6111
$ws.Range("K107").Value = 1307.6111
$ws.Range("M107").Value = 612.3888999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 243093.81
$ws.Range("I31").Value = 384472.88
$ws.Range("J31").Value = 3364.087
$ws.Range("K31").Value = 384472.88
$ws.Range("L31").Value = 3364.087
$ws.Range("M31").Value = -384177.88
$ws.Range("N31").Value = -3954.087
$ws.Range("H34").Value = 243093.81
$ws.Range("I34").Value = 384472.88
$ws.Range("J34").Value = 3364.087
$ws.Range("K34").Value = 384472.88
$ws.Range("L34").Value = 3364.087
$ws.Range("M34").Value = -384270.88
$ws.Range("N34").Value = -3768.087
$ws.Range("H58").Value = 2952.1707
$ws.Range("I58").Value = 2287.0625
$ws.Range("J58").Value = 3377.84
$ws.Range("K58").Value = 2287.0625
$ws.Range("L58").Value = 3377.84
$ws.Range("M58").Value = -2084.0625
$ws.Range("N58").Value = -3783.84
$ws.Range("H99").Value = 18388.785
$ws.Range("I99").Value = 35357.5
$ws.Range("K99").Value = 35357.5
$ws.Range("M99").Value = -33859.5
$ws.Range("H109").Value = 34999.5
$ws.Range("J109").Value = 34999.5
$ws.Range("L109").Value = 34999.5
$ws.Range("N109").Value = -37079.5
$ws.Range("H126").Value = 18388.785
$ws.Range("I126").Value = 35357.5
$ws.Range("K126").Value = 106072.5
$ws.Range("M126").Value = -103602.5
$ws.Range("H134").Value = 3490
$ws.Range("I134").Value = 2177.44
$ws.Range("K134").Value = 6532.32
$ws.Range("M134").Value = -3997.32
$ws.Range("H136").Value = 2952.1707
$ws.Range("I136").Value = 2287.0625
$ws.Range("J136").Value = 3377.84
$ws.Range("K136").Value = 6861.1875
$ws.Range("L136").Value = 10133.52
$ws.Range("M136").Value = -4311.1875
$ws.Range("N136").Value = -15233.52

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 192.92308
$ws.Range("I14").Value = 192.92308
$ws.Range("K14").Value = 578.76924
$ws.Range("M14").Value = -405.76924
$ws.Range("H62").Value = 2415.8333
$ws.Range("I62").Value = 1599.4
$ws.Range("K62").Value = 4798.200000000001
$ws.Range("M62").Value = -4112.200000000001
$ws.Range("H65").Value = 2415.8333
$ws.Range("I65").Value = 1599.4
$ws.Range("K65").Value = 14394.6
$ws.Range("M65").Value = -10962.6
$ws.Range("H80").Value = 10000
$ws.Range("J80").Value = 10000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31872
$ws.Range("H83").Value = 10000
$ws.Range("J83").Value = 10000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99360
$ws.Range("H125").Value = 6015
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H128").Value = 254711
$ws.Range("I128").Value = 254711
$ws.Range("K128").Value = 764133
$ws.Range("M128").Value = -759153
$ws.Range("H133").Value = 9813.299999999999
$ws.Range("I133").Value = 3159.25
$ws.Range("J133").Value = 14249.333
$ws.Range("K133").Value = 9477.75
$ws.Range("L133").Value = 42747.999
$ws.Range("M133").Value = -4417.75
$ws.Range("N133").Value = -52867.999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19528.889
$ws.Range("I70").Value = 22290.2
$ws.Range("J70").Value = 5722.3335
$ws.Range("K70").Value = 22290.2
$ws.Range("L70").Value = 5722.3335
$ws.Range("M70").Value = -22020.2
$ws.Range("N70").Value = -6262.3335
$ws.Range("H73").Value = 19528.889
$ws.Range("I73").Value = 22290.2
$ws.Range("J73").Value = 5722.3335
$ws.Range("K73").Value = 22290.2
$ws.Range("L73").Value = 5722.3335
$ws.Range("M73").Value = -21354.2
$ws.Range("N73").Value = -7594.3335
$ws.Range("H126").Value = 3388.4443
$ws.Range("J126").Value = 3249.75
$ws.Range("L126").Value = 9749.25
$ws.Range("N126").Value = -14689.25
$ws.Range("H132").Value = 11343.725
$ws.Range("I132").Value = 9136.120000000001
$ws.Range("J132").Value = 15023.066
$ws.Range("K132").Value = 27408.36
$ws.Range("L132").Value = 45069.198
$ws.Range("M132").Value = -24878.36
$ws.Range("N132").Value = -50129.198

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 3484.3333
$ws.Range("J55").Value = 3484.3333
$ws.Range("L55").Value = 3484.3333
$ws.Range("N55").Value = -4038.3333
$ws.Range("H126").Value = 2551.8572
$ws.Range("I126").Value = 1566.35
$ws.Range("J126").Value = 5015.625
$ws.Range("K126").Value = 4699.049999999999
$ws.Range("L126").Value = 15046.875
$ws.Range("M126").Value = -2229.049999999999
$ws.Range("N126").Value = -19986.875
$ws.Range("H136").Value = 55558476
$ws.Range("I136").Value = 76925624
$ws.Range("J136").Value = 3896.8
$ws.Range("K136").Value = 230776872
$ws.Range("L136").Value = 11690.4
$ws.Range("M136").Value = -230774322
$ws.Range("N136").Value = -16790.4

Write-Host "Applied 222 cell updates across 7 sheets"
